# maj planning + script arthur
#
# Planning.xlsx update:
#  - "Scoring" (A38) renamed to "Scoring + timer"
#  - "Mécanisme de retry apres un lancé" (row 32) marked done: B32 -> 100% (green),
#    assignee note D32 cleared
#  - "Score" (row 46) marked done: B46 -> 100% (green), assignee note D46 cleared
#  - selection moved to A38 (scrolled view around row 22)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- A38: "Scoring" -> "Scoring + timer" ---
$ws.Range("A38").Value = "Scoring + timer"

# --- Row 32 ("Mécanisme de retry apres un lancé") completed ---
# Copy the "done" (green, 100%) percentage format from another completed row
# so B32 picks up the same cell style Excel uses elsewhere (s=4), then set it to 100%.
$ws.Range("B9").Copy()
$ws.Range("B32").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B32").Value = 1
$ws.Range("D32").Clear()               # assignee no longer needed once task is done

# --- Row 46 ("Score") completed ---
$ws.Range("B9").Copy()
$ws.Range("B46").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B46").Value = 1
$ws.Range("D46").Clear()               # assignee no longer needed once task is done

$excel.CutCopyMode = 0

# --- Update the view: scroll near row 22, select A38 ---
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A38").Select()
